$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KEYWORDS (column B) and ACTIONS (column C) for rule R1 (row 2)
$ws.Range("B2").Value = "[LDAP: error code 49 - 80090308: LdapErr: DSID-0C09042F, comment: AcceptSecurityContext error, data 531, v2580 ]"
$ws.Range("C2").Value = "1. Ask client for microsoft error debugging."

# Update KEYWORDS (column B) and ACTIONS (column C) for rule R2 (row 3)
$ws.Range("B3").Value = "Security token is invalid. java.util.NoSuchElementException: No value present"
$ws.Range("C3").Value = "1. Ask client for Security token ."

# Move the active selection to C2 (was C4)
$ws.Range("C2").Select()
